$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C, matching the style used by A1/B1
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Coord: normal vector scan"

# Updated "Angle between normals (degrees)" values (column B)
$ws.Range("B2").Value = 0.3475469908725779
$ws.Range("B3").Value = 1.204601251540679
$ws.Range("B4").Value = 0.218795725544832
$ws.Range("B5").Value = 0.6007661531488215
$ws.Range("B6").Value = 0.926475704553728
$ws.Range("B7").Value = 0.7115700705739271
$ws.Range("B8").Value = 0.4241854348423081
$ws.Range("B9").Value = 1.137964285050551
$ws.Range("B10").Value = 0.4038639736697666
$ws.Range("B11").Value = 1.003521588268715
$ws.Range("B12").Value = 0.5989561537111915
$ws.Range("B13").Value = 1.812873385486746
$ws.Range("B14").Value = 3.472969097681784
$ws.Range("B15").Value = 1.265110129265497

# New "Coord: normal vector scan" values (column C)
$ws.Range("C2").Value = "[0.         0.31922366 0.9476794 ]"
$ws.Range("C3").Value = "[-0.44173296  0.5244037   0.72792359]"
$ws.Range("C4").Value = "[-0.00312578  0.01384813  0.99989922]"
$ws.Range("C5").Value = "[-6.33261181e-04  2.84442762e-01 -9.58692815e-01]"
$ws.Range("C6").Value = "[0.73091311 0.2947569  0.61553586]"
$ws.Range("C7").Value = "[-0.74681708 -0.28082557  0.60282771]"
$ws.Range("C8").Value = "[0.         0.31795577 0.94810555]"
$ws.Range("C9").Value = "[ 0.         -0.30612011  0.9519929 ]"
$ws.Range("C10").Value = "[-0.73998515  0.26934801  0.61633889]"
$ws.Range("C11").Value = "[ 0.73531595 -0.28361354  0.61552726]"
$ws.Range("C12").Value = "[ 6.30353532e-04 -2.84472935e-01 -9.58683865e-01]"
$ws.Range("C13").Value = "[-0.72311312 -0.28745969  0.62807193]"
$ws.Range("C14").Value = "[0.70523896 0.26133167 0.65904763]"
$ws.Range("C15").Value = "[ 0.         -0.30400678  0.95266987]"
